$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Transacciones")

# --- Add new transaction row (row 34) to "Transacciones" sheet ---
# Copy formatting (date style) from the row above, then set values.
$ws1.Range("A33").Copy($ws1.Range("A34")) | Out-Null

$ws1.Range("A34").Value = 43566
$ws1.Range("B34").Value = 38
$ws1.Range("C34").Value = "Garrafón Ciel"
$ws1.Range("D34").Value = "Despensa"
$ws1.Range("E34").Value = "Gasto"
$ws1.Range("F34").Value = "Efectivo"
$ws1.Range("G34").Value = "Extra"
$ws1.Range("K34").Value = 7358.64
$ws1.Range("L34").Value = 61
$ws1.Range("M34").Formula = "=M33-B34"
$ws1.Range("N34").Formula = "=SUM(K34:M34)"
$ws1.Range("O34").Formula = "=N34-4000"

# --- Switch the active sheet from "Deudas" back to "Transacciones" ---
$ws1.Activate() | Out-Null
$ws1.Range("M35").Select() | Out-Null
